# Apply the diff:
#  1. Insert a new "Player Info" sheet in front of the existing sheets,
#     with player bio columns ID / NAME / BATTING_HAND / BOWL_STYLE.
#  2. On the existing "ODI Batting" and "ODI Bowling" sheets, rename the
#     MATCH_CARD_LINK column to MATCH_CODE and replace each full scorecard
#     URL with just the numeric match code pulled from its
#     "MatchCode=" query parameter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "Player Info" sheet before the current first sheet
#    ("ODI Batting"), so the final tab order is:
#      Player Info, ODI Batting, ODI Bowling
# ---------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

# Header row
$playerInfo.Cells.Item(1, 1).Value = "ID"
$playerInfo.Cells.Item(1, 2).Value = "NAME"
$playerInfo.Cells.Item(1, 3).Value = "BATTING_HAND"
$playerInfo.Cells.Item(1, 4).Value = "BOWL_STYLE"

# Match the bold / bordered / centered look used by the header row on the
# other sheets.
$playerInfoHeader = $playerInfo.Range("A1:D1")
$playerInfoHeader.Font.Bold = $true
$playerInfoHeader.Borders.LineStyle = 1
$playerInfoHeader.HorizontalAlignment = -4108
$playerInfoHeader.VerticalAlignment = -4160

# Data row. ID is written as text (like the numeric-looking text columns,
# e.g. MATCH_NUMBER, on the other sheets): mark the cell as Text first so
# "3829" isn't auto-converted to a number, then drop back to the default
# "Normal" style so no stray number-format style sticks to the cell.
$idCell = $playerInfo.Cells.Item(2, 1)
$idCell.NumberFormat = "@"
$idCell.Value = "3829"
$idCell.Style = "Normal"

$playerInfo.Cells.Item(2, 2).Value = "Hamish Kyle Bennett"
$playerInfo.Cells.Item(2, 3).Value = "Left Handed"
$playerInfo.Cells.Item(2, 4).Value = "Right Arm Medium Fast"

# ---------------------------------------------------------------------
# 2. "ODI Batting" sheet (now the 2nd tab): MATCH_CARD_LINK (col D) ->
#    MATCH_CODE, values become just the MatchCode number.
# ---------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Cells.Item(1, 4).Value = "MATCH_CODE"

$battingLastRow = $battingSheet.UsedRange.Rows.Count
$battingCodeRange = $battingSheet.Range("D2:D" + $battingLastRow)
$battingCodeRange.NumberFormat = "@"
for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $link = $cell.Text
    if ($link -and $link -like "*MatchCode=*") {
        $code = ($link -split "MatchCode=")[1]
        $cell.Value = $code
    }
}
$battingCodeRange.Style = "Normal"

# ---------------------------------------------------------------------
# 3. "ODI Bowling" sheet (now the 3rd tab): MATCH_CARD_LINK (col B) ->
#    MATCH_CODE, values become just the MatchCode number.
# ---------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Cells.Item(1, 2).Value = "MATCH_CODE"

$bowlingLastRow = $bowlingSheet.UsedRange.Rows.Count
$bowlingCodeRange = $bowlingSheet.Range("B2:B" + $bowlingLastRow)
$bowlingCodeRange.NumberFormat = "@"
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $link = $cell.Text
    if ($link -and $link -like "*MatchCode=*") {
        $code = ($link -split "MatchCode=")[1]
        $cell.Value = $code
    }
}
$bowlingCodeRange.Style = "Normal"

# Leave selection on A1 of the first sheet, consistent with the original file
$playerInfo.Range("A1").Select()
